$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.456404
$ws.Range("H2").Value = 1.369212
$ws.Range("I2").Value = 0.01914960767004715
$ws.Range("J2").Value = 0.01914960767004715
$ws.Range("M2").Value = 0.668273
$ws.Range("N2").Value = 2.004819
$ws.Range("O2").Value = 0.01328414746766746
$ws.Range("P2").Value = 0.01328414746766746
$ws.Range("Q2").Value = 0.305002470292
$ws.Range("R2").Value = 2.745022232628
$ws.Range("S2").Value = 0.0002543862122368823
$ws.Range("T2").Value = 0.0002543862122368822

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.456404
$ws.Range("H3").Value = 1.369212
$ws.Range("I3").Value = 0.01914960767004715
$ws.Range("J3").Value = 0.01914960767004715
$ws.Range("O3").Value = 0.3831531055114357
$ws.Range("P3").Value = 0.3831531055114357
$ws.Range("Q3").Value = 8.797150435545333
$ws.Range("R3").Value = 79.174353919908
$ws.Range("S3").Value = 0.007337231648104175
$ws.Range("T3").Value = 0.007337231648104174

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.456404
$ws.Range("H4").Value = 1.369212
$ws.Range("I4").Value = 0.01914960767004715
$ws.Range("J4").Value = 0.01914960767004715
$ws.Range("M4").Value = 30.36285833333334
$ws.Range("N4").Value = 91.08857500000001
$ws.Range("O4").Value = 0.6035627470208969
$ws.Range("P4").Value = 0.6035627470208967
$ws.Range("Q4").Value = 13.85772999476667
$ws.Range("R4").Value = 124.7195699529
$ws.Range("S4").Value = 0.01155798980970609
$ws.Range("T4").Value = 0.01155798980970609

# Row 5
$ws.Range("I5").Value = 0.8285024587002443
$ws.Range("J5").Value = 0.8285024587002443
$ws.Range("M5").Value = 0.668273
$ws.Range("N5").Value = 2.004819
$ws.Range("O5").Value = 0.01328414746766746
$ws.Range("P5").Value = 0.01328414746766746
$ws.Range("Q5").Value = 13.19584718917367
$ws.Range("R5").Value = 118.762624702563
$ws.Range("S5").Value = 0.01100594883869912
$ws.Range("T5").Value = 0.01100594883869911

# Row 6
$ws.Range("I6").Value = 0.8285024587002443
$ws.Range("J6").Value = 0.8285024587002443
$ws.Range("O6").Value = 0.3831531055114357
$ws.Range("P6").Value = 0.3831531055114357
$ws.Range("S6").Value = 0.3174432899748586
$ws.Range("T6").Value = 0.3174432899748586

# Row 7
$ws.Range("I7").Value = 0.8285024587002443
$ws.Range("J7").Value = 0.8285024587002443
$ws.Range("M7").Value = 30.36285833333334
$ws.Range("N7").Value = 91.08857500000001
$ws.Range("O7").Value = 0.6035627470208969
$ws.Range("P7").Value = 0.6035627470208967
$ws.Range("Q7").Value = 599.550840439753
$ws.Range("R7").Value = 5395.957563957776
$ws.Range("S7").Value = 0.5000532198866866
$ws.Range("T7").Value = 0.5000532198866865

# Row 8
$ws.Range("I8").Value = 0.1523479336297086
$ws.Range("J8").Value = 0.1523479336297086
$ws.Range("M8").Value = 0.668273
$ws.Range("N8").Value = 2.004819
$ws.Range("O8").Value = 0.01328414746766746
$ws.Range("P8").Value = 0.01328414746766746
$ws.Range("Q8").Value = 2.426498594727001
$ws.Range("R8").Value = 21.838487352543
$ws.Range("S8").Value = 0.002023812416731464
$ws.Range("T8").Value = 0.002023812416731464

# Row 9
$ws.Range("I9").Value = 0.1523479336297086
$ws.Range("J9").Value = 0.1523479336297086
$ws.Range("O9").Value = 0.3831531055114357
$ws.Range("P9").Value = 0.3831531055114357
$ws.Range("S9").Value = 0.05837258388847296
$ws.Range("T9").Value = 0.05837258388847295

# Row 10
$ws.Range("I10").Value = 0.1523479336297086
$ws.Range("J10").Value = 0.1523479336297086
$ws.Range("M10").Value = 30.36285833333334
$ws.Range("N10").Value = 91.08857500000001
$ws.Range("O10").Value = 0.6035627470208969
$ws.Range("P10").Value = 0.6035627470208967
$ws.Range("R10").Value = 992.2275742092752
$ws.Range("S10").Value = 0.09195153732450422
$ws.Range("T10").Value = 0.0919515373245042
